$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.175.96'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.749.89'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.39%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '602.09'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.81'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.19%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.749.01'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.42%  '
$ws.Range('E9').Value = '  +1.21%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.169'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.53%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.38'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.33%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.459'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.48%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '38.02'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.64%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000248'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.72%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.381.54'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.31%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.748.76'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.11%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '69.187.87'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.27%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.37'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.74%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.46'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.39%  '
$ws.Range('E20').Value = '  -1.47%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.14'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +18.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '494.38'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.51%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.727'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.33%  '
$ws.Range('E24').Value = '  +7.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.87'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('E26').Value = '  -0.56%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.30'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.10'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('E30').Value = '  +1.86%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.47'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.39%  '
$ws.Range('E32').Value = '  +1.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '31.60'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.09%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.897.29'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.690.18'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.20%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.108'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.57%  '
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.98'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.24%  '
$ws.Range('E39').Value = '  +0.79%  '
$ws.Range('E40').Value = '  +1.06%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.324'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('E42').Value = '  +4.51%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '48.82'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.52%  '
$ws.Range('B44').Value = 'Bittensor'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '429.00'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.29%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.98'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.39%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.48'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.17%  '
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '40.21'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.48%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '141.17'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.43%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.795.09'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.47%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0352'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.06%  '
